$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)
$ws.Activate()

# New test case block appended below the existing tc_06 data (rows 1-9)
$ws.Range("A11").Value = "verifyLoginWithValidCred"
$ws.Range("A12").Value = "username"
$ws.Range("B12").Value = "password"
$ws.Range("A13").Value = "w2ajava@way2automation.com"
$ws.Range("B13").Value = "Tcs@12345"

# Turn the credential cell into a mailto hyperlink (adds Hyperlink style/font automatically)
$ws.Hyperlinks.Add($ws.Range("B13"), "mailto:w2ajava@way2automation.com")

# Match the recorded selection state from the saved workbook
$ws.Range("B13").Select() | Out-Null
